$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Dropout: A Simple Way to Prevent Neural Networks from`nOverfitting"
$ws.Range("B5").Value = "用来防止过拟合"

$ws.Range("A5").WrapText = $true

$ws.Rows.Item(2).RowHeight = 28.5
$ws.Rows.Item(5).RowHeight = 28.5

$ws.Range("K9").Select()
